$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string content changes -----------------------------------

# Swap the Sri Lanka / Venezuela rows (A100 <-> A101 text) so the country
# list order becomes: ... Cuba, Venezuela, Sri Lanka, Afganistan ...
$ws.Range("A100").Value = "Venezuela"
$ws.Range("A101").Value = "Sri Lanka"

# Update the "last updated" timestamp string (22:29 -> 22:59).
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 22:59"

# --- Data refresh (updated case counts) -------------------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 120529
$ws.Range("C4").Value = 16403
$ws.Range("E4").Value = 115292
$ws.Range("G4").Value = 312
$ws.Range("H4").Value = 2008

# Row 7: España
$ws.Range("B7").Value = 72469
$ws.Range("C7").Value = 6750
$ws.Range("E7").Value = 54358
$ws.Range("G7").Value = 688
$ws.Range("H7").Value = 5826

# Row 18: Canada
$ws.Range("E18").Value = 5162
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 60

# Row 20: Noruega
$ws.Range("B20").Value = 4013
$ws.Range("C20").Value = 242
$ws.Range("E20").Value = 3983

# Row 39: Sudafrica
$ws.Range("E39").Value = 1155
$ws.Range("H39").Value = 1

# Row 100: now Venezuela (updated figures)
$ws.Range("B100").Value = 119
$ws.Range("C100").Value = 6
$ws.Range("D100").Value = 39
$ws.Range("E100").Value = 78
$ws.Range("F100").Value = 2
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 2

# Row 101: now Sri Lanka (updated figures)
$ws.Range("B101").Value = 113
$ws.Range("C101").Value = 7
$ws.Range("D101").Value = 9
$ws.Range("E101").Value = 103
$ws.Range("F101").Value = 5
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 1

# Row 108: Estado de Palestina
$ws.Range("B108").Value = 98
$ws.Range("C108").Value = 7
$ws.Range("E108").Value = 79
